$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.489.48'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').Value = '3.355.05'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'400.51"
$ws.Range('E5').Value = '  -3.52%  '
$ws.Range('D6').Value = "'126.03"
$ws.Range('E6').Value = '  +7.67%  '
$ws.Range('E7').Value = '  +2.06%  '
$ws.Range('D8').Value = "'0.999"
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +4.49%  '
$ws.Range('E10').Value = '  +1.78%  '
$ws.Range('D11').Value = "'40.91"
$ws.Range('E11').Value = '  +2.17%  '
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('D13').Value = '3.881.72'
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('E14').Value = '  -1.25%  '
$ws.Range('E15').Value = '  -0.46%  '
$ws.Range('D16').Value = '3.351.86'
$ws.Range('E16').Value = '  -1.17%  '
$ws.Range('D17').Value = '61.430.54'
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('E18').Value = '  +2.64%  '
$ws.Range('E19').Value = '  -0.69%  '
$ws.Range('E20').Value = '  +7.22%  '
$ws.Range('E21').Value = '  -4.63%  '
$ws.Range('D22').Value = "'79.86"
$ws.Range('E22').Value = '  +6.62%  '
$ws.Range('D23').Value = "'12.70"
$ws.Range('E23').Value = '  +0.65%  '
$ws.Range('D24').Value = "'298.54"
$ws.Range('E24').Value = '  +0.32%  '
$ws.Range('D25').Value = "'3.09"
$ws.Range('E25').Value = '  -1.71%  '
$ws.Range('E26').Value = '  +11.25%  '
$ws.Range('B27').Value = 'Filecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D27').Value = "'8.20"
$ws.Range('E27').Value = '  +7.69%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'28.97"
$ws.Range('E28').Value = '  -2.40%  '
$ws.Range('D29').Value = "'7.45"
$ws.Range('E29').Value = '  -6.35%  '
$ws.Range('E30').Value = '  -3.13%  '
$ws.Range('E31').Value = '  +1.06%  '
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('D33').Value = "'11.30"
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('D34').Value = "'2.50"
$ws.Range('E34').Value = '  -2.08%  '
$ws.Range('D35').Value = "'40.98"
$ws.Range('E35').Value = '  -5.17%  '
$ws.Range('D36').Value = "'0.0478"
$ws.Range('E36').Value = '  -2.65%  '
$ws.Range('D37').Value = "'51.96"
$ws.Range('E37').Value = '  -0.61%  '
$ws.Range('D38').Value = "'0.999"
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('D39').Value = "'3.37"
$ws.Range('E39').Value = '  -2.17%  '
$ws.Range('E40').Value = '  -6.49%  '
$ws.Range('D41').Value = "'137.05"
$ws.Range('E41').Value = '  +2.64%  '
$ws.Range('E42').Value = '  +2.23%  '
$ws.Range('E43').Value = '  +0.97%  '
$ws.Range('E44').Value = '  -2.64%  '
$ws.Range('D45').Value = "'3.89"
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('D46').Value = "'16.56"
$ws.Range('E46').Value = '  +0.35%  '
$ws.Range('E47').Value = '  -0.49%  '
$ws.Range('D48').Value = "'20.93"
$ws.Range('E48').Value = '  -1.56%  '
$ws.Range('D49').Value = '3.684.45'
$ws.Range('E49').Value = '  -0.39%  '
$ws.Range('D50').Value = '2.095.74'
$ws.Range('E50').Value = '  -3.52%  '
$ws.Range('E51').Value = '  -4.77%  '
